# Updated cryptos list on Thu Feb 15 05:54:54 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns of the crypto ranking
# table, and fixes rows 5/6 where Solana and BNB had been swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.941.78"
$ws.Range("E2").Value = "  +4.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.773.32"
$ws.Range("E3").Value = "  +5.08%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "340.04"
$ws.Range("E5").Value = "  +4.29%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.43"
$ws.Range("E6").Value = "  +2.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.546"
$ws.Range("E7").Value = "  +4.49%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  +4.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.59"
$ws.Range("E10").Value = "  +4.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("E11").Value = "  +5.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.95"
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("E13").Value = "  +2.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.56"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.216.10"
$ws.Range("E15").Value = "  +5.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.792.94"
$ws.Range("E16").Value = "  +5.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "51.791.89"
$ws.Range("E17").Value = "  +4.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.876"
$ws.Range("E18").Value = "  +2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.17"
$ws.Range("E19").Value = "  +8.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.20"
$ws.Range("E20").Value = "  -1.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.93"
$ws.Range("E21").Value = "  +3.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  +2.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "275.91"
$ws.Range("E23").Value = "  +2.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.82"
$ws.Range("E24").Value = "  +1.04%  "

$ws.Range("E25").Value = "  +6.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.58"
$ws.Range("E26").Value = "  +2.06%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("E29").Value = "  +0.98%  "

$ws.Range("E30").Value = "  +2.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.57"
$ws.Range("E31").Value = "  -0.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.03"
$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.68"
$ws.Range("E33").Value = "  +3.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0821"
$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.09"
$ws.Range("E36").Value = "  +3.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.78"
$ws.Range("E37").Value = "  -1.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.91"
$ws.Range("E38").Value = "  -0.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.21"
$ws.Range("E39").Value = "  +3.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0375"
$ws.Range("E40").Value = "  +9.95%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("E41").Value = "  +24.24%  "

$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("E43").Value = "  +3.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.58"
$ws.Range("E44").Value = "  -3.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.98"
$ws.Range("E45").Value = "  -0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.072.63"
$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("E48").Value = "  +3.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.51"
$ws.Range("E49").Value = "  +5.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.83"
$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.02"
$ws.Range("E51").Value = "  +0.51%  "
